$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 448, shifting existing rows 448:486 down to 449:487
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with the new weekly record.
# Columns A,B,C,E,F,G,H,N,O,Q,R are constant across this block of rows.
$ws.Range("A448").Value = 4
$ws.Range("B448").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C448").Value = "Los Lagos"
$ws.Range("D448").Value = 45106
$ws.Range("E448").Value = 10
$ws.Range("F448").Value = 100112017
$ws.Range("G448").Value = "Apio"
$ws.Range("H448").Value = "Americana (o)"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 25
$ws.Range("K448").Value = 11000
$ws.Range("L448").Value = 11000
$ws.Range("M448").Value = 11000
$ws.Range("N448").Value = "$/docena de matas"
$ws.Range("O448").Value = "Región de Coquimbo"
$ws.Range("P448").Value = 1833
$ws.Range("Q448").Value = 6
$ws.Range("R448").Value = "Hortaliza"
